# Reorders the word list in column A (rows 16-56) to match the
# re-exported shared-strings order described in the commit diff.
# Column B (the counts) is left untouched; only the word labels that
# sit in tied-count groups are reshuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newWords = @(
    "колеса",
    "Крымскую соль",
    "полотно",
    "парча",
    "говядина",
    "сено",
    "табак",
    "позумент",
    "выбойка",
    "чулок",
    "сахар",
    "шелк",
    "лыко",
    "лес",
    "сапог",
    "сани",
    "ладан",
    "коса",
    "китайка",
    "замок",
    "овца",
    "веревка",
    "конь",
    "горшок",
    "обод",
    "платок",
    "рогожа",
    "гвоздь",
    "ром",
    "сосуд",
    "бечева",
    "хомут",
    "брусья",
    "гумми",
    "скотский кожа",
    "котел",
    "покроми",
    "роза",
    "нитка",
    "сковорода",
    "дуга"
)

$startRow = 16
for ($i = 0; $i -lt $newWords.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newWords[$i]
}
